$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows (12-24) that are no longer part of the data set
$ws.Rows("12:24").Delete()

# Write the new policy id values into A2:A11
$values = @(
    "100-0000043",
    "100-0000044",
    "100-0000045",
    "100-0000046",
    "100-0000047",
    "100-0000048",
    "100-0000049",
    "100-0000050",
    "100-0000051",
    "100-0000052"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Reset the view: scroll back to the top-left and move the selection to C7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
